# Fruta / hortaliza, semanal
# Insert two new weekly price records at row 252 (pushing the existing
# rows 252-255 down to rows 254-257), then populate the two new rows
# with the new "Primera"/"Segunda" quality records dated 44656.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 252; this shifts the old
# rows 252-255 down to 254-257 and keeps their values/styles intact.
$ws.Range("A252:A253").EntireRow.Insert()

# New row 252: Primera
$ws.Cells.Item(252, 1).Value = 4
$ws.Cells.Item(252, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(252, 3).Value = "Los Lagos"
$ws.Cells.Item(252, 4).Value = 44656
$ws.Cells.Item(252, 5).Value = 10
$ws.Cells.Item(252, 6).Value = "Fruta"
$ws.Cells.Item(252, 7).Value = 100102
$ws.Cells.Item(252, 8).Value = "Cítricos"
$ws.Cells.Item(252, 9).Value = 100102006
$ws.Cells.Item(252, 10).Value = "Pomelo"
$ws.Cells.Item(252, 11).Value = "Start Ruby"
$ws.Cells.Item(252, 12).Value = "Primera"
$ws.Cells.Item(252, 13).Value = 200
$ws.Cells.Item(252, 14).Value = 14000
$ws.Cells.Item(252, 15).Value = 15000
$ws.Cells.Item(252, 16).Value = 14500
$ws.Cells.Item(252, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(252, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(252, 19).Value = 1036
$ws.Cells.Item(252, 20).Value = 14

# New row 253: Segunda
$ws.Cells.Item(253, 1).Value = 4
$ws.Cells.Item(253, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(253, 3).Value = "Los Lagos"
$ws.Cells.Item(253, 4).Value = 44656
$ws.Cells.Item(253, 5).Value = 10
$ws.Cells.Item(253, 6).Value = "Fruta"
$ws.Cells.Item(253, 7).Value = 100102
$ws.Cells.Item(253, 8).Value = "Cítricos"
$ws.Cells.Item(253, 9).Value = 100102006
$ws.Cells.Item(253, 10).Value = "Pomelo"
$ws.Cells.Item(253, 11).Value = "Start Ruby"
$ws.Cells.Item(253, 12).Value = "Segunda"
$ws.Cells.Item(253, 13).Value = 100
$ws.Cells.Item(253, 14).Value = 12000
$ws.Cells.Item(253, 15).Value = 12000
$ws.Cells.Item(253, 16).Value = 12000
$ws.Cells.Item(253, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(253, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(253, 19).Value = 857
$ws.Cells.Item(253, 20).Value = 14
